$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 450.8
$ws.Range("I2").Value = 338.5
$ws.Range("K2").Value = 338.5
$ws.Range("M2").Value = -225.5
$ws.Range("H19").Value = 10948.546
$ws.Range("I19").Value = 1941.2858
$ws.Range("J19").Value = 15151.934
$ws.Range("K19").Value = 1941.2858
$ws.Range("L19").Value = 15151.934
$ws.Range("M19").Value = -1766.2858
$ws.Range("N19").Value = -15501.934
$ws.Range("H70").Value = 3402198
$ws.Range("I70").Value = 3402198
$ws.Range("K70").Value = 10206594
$ws.Range("M70").Value = -10206324
$ws.Range("H73").Value = 3402198
$ws.Range("I73").Value = 3402198
$ws.Range("K73").Value = 10206594
$ws.Range("M73").Value = -10205658
$ws.Range("H132").Value = 4361.8667
$ws.Range("I132").Value = 4025.52
$ws.Range("J132").Value = 6043.6
$ws.Range("K132").Value = 12076.56
$ws.Range("L132").Value = 18130.8
$ws.Range("M132").Value = -9546.559999999999
$ws.Range("N132").Value = -23190.8
$ws.Range("H137").Value = 52632676
$ws.Range("I137").Value = 52632676
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 157898028
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -157895478
$ws.Range("N137").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 25642774
$ws.Range("I61").Value = 25642774
$ws.Range("K61").Value = 25642774
$ws.Range("M61").Value = -25642562
$ws.Range("H63").Value = 85722296
$ws.Range("I63").Value = 111113200
$ws.Range("J63").Value = 40018676
$ws.Range("K63").Value = 111113200
$ws.Range("L63").Value = 40018676
$ws.Range("M63").Value = -111112514
$ws.Range("N63").Value = -40020048
$ws.Range("H66").Value = 85722296
$ws.Range("I66").Value = 111113200
$ws.Range("J66").Value = 40018676
$ws.Range("K66").Value = 555566000
$ws.Range("L66").Value = 200093380
$ws.Range("M66").Value = -555562568
$ws.Range("N66").Value = -200100244
$ws.Range("H74").Value = 2498.724
$ws.Range("I74").Value = 1664.5454
$ws.Range("J74").Value = 5120.4287
$ws.Range("K74").Value = 1664.5454
$ws.Range("L74").Value = 5120.4287
$ws.Range("M74").Value = -790.5454
$ws.Range("N74").Value = -6868.4287
$ws.Range("H77").Value = 2498.724
$ws.Range("I77").Value = 1664.5454
$ws.Range("J77").Value = 5120.4287
$ws.Range("K77").Value = 8322.726999999999
$ws.Range("L77").Value = 25602.1435
$ws.Range("M77").Value = -3954.726999999999
$ws.Range("N77").Value = -34338.14350000001
$ws.Range("H97").Value = 723
$ws.Range("I97").Value = 519.36365
$ws.Range("K97").Value = 519.36365
$ws.Range("M97").Value = -23.36365000000001
$ws.Range("H110").Value = 68830.2
$ws.Range("I110").Value = 43787.832
$ws.Range("J110").Value = 168999.67
$ws.Range("K110").Value = 43787.832
$ws.Range("L110").Value = 168999.67
$ws.Range("M110").Value = -41742.832
$ws.Range("N110").Value = -173089.67
$ws.Range("H136").Value = 25642774
$ws.Range("I136").Value = 25642774
$ws.Range("K136").Value = 76928322
$ws.Range("M136").Value = -76925772

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 2540
$ws.Range("L80").Value = 2540
$ws.Range("N80").Value = -4536
$ws.Range("H81").Value = 56998
$ws.Range("J81").Value = 56998
$ws.Range("L81").Value = 56998
$ws.Range("N81").Value = -59120
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 2540
$ws.Range("L83").Value = 12700
$ws.Range("N83").Value = -22684
$ws.Range("H84").Value = 56998
$ws.Range("J84").Value = 56998
$ws.Range("L84").Value = 170994
$ws.Range("N84").Value = -181602
$ws.Range("H94").Value = 5382.2856
$ws.Range("I94").Value = 5953.2
$ws.Range("K94").Value = 5953.2
$ws.Range("M94").Value = -5502.2
$ws.Range("H107").Value = 25264.428
$ws.Range("I107").Value = 18420.176
$ws.Range("K107").Value = 18420.176
$ws.Range("M107").Value = -16500.176
$ws.Range("H128").Value = 5000
$ws.Range("I128").Value = 5000
$ws.Range("K128").Value = 15000
$ws.Range("M128").Value = -12510
$ws.Range("H134").Value = 1689.2444
$ws.Range("I134").Value = 1448.325
$ws.Range("K134").Value = 4344.975
$ws.Range("M134").Value = -1809.975

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3095.1482
$ws.Range("I31").Value = 2361.875
$ws.Range("J31").Value = 3403.8948
$ws.Range("K31").Value = 2361.875
$ws.Range("L31").Value = 3403.8948
$ws.Range("M31").Value = -2066.875
$ws.Range("N31").Value = -3993.8948
$ws.Range("H34").Value = 3095.1482
$ws.Range("I34").Value = 2361.875
$ws.Range("J34").Value = 3403.8948
$ws.Range("K34").Value = 2361.875
$ws.Range("L34").Value = 3403.8948
$ws.Range("M34").Value = -2159.875
$ws.Range("N34").Value = -3807.8948
$ws.Range("H59").Value = 54999.5
$ws.Range("I59").Value = 54999.5
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 54999.5
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -53854.5
$ws.Range("N59").ClearContents()
$ws.Range("H62").Value = 6092
$ws.Range("J62").Value = 4765.6
$ws.Range("L62").Value = 4765.6
$ws.Range("N62").Value = -6013.6
$ws.Range("H65").Value = 6092
$ws.Range("J65").Value = 4765.6
$ws.Range("L65").Value = 23828
$ws.Range("N65").Value = -30068
$ws.Range("H107").Value = 1676.2667
$ws.Range("I107").Value = 1557.3158
$ws.Range("J107").Value = 1881.7273
$ws.Range("K107").Value = 1557.3158
$ws.Range("L107").Value = 1881.7273
$ws.Range("M107").Value = 362.6841999999999
$ws.Range("N107").Value = -5721.7273
$ws.Range("H132").Value = 2829.3635
$ws.Range("I132").Value = 2726.0476
$ws.Range("K132").Value = 8178.1428
$ws.Range("M132").Value = -5648.1428

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 271.77777
$ws.Range("I92").Value = 223.33333
$ws.Range("J92").Value = 368.66666
$ws.Range("K92").Value = 669.99999
$ws.Range("L92").Value = 1105.99998
$ws.Range("M92").Value = 578.00001
$ws.Range("N92").Value = -3601.99998
$ws.Range("H98").Value = 7570.125
$ws.Range("I98").Value = 11582.6
$ws.Range("J98").Value = 882.6667
$ws.Range("K98").Value = 34747.8
$ws.Range("L98").Value = 2648.0001
$ws.Range("M98").Value = -33249.8
$ws.Range("N98").Value = -5644.0001
$ws.Range("H132").Value = 1475
$ws.Range("I132").Value = 1475
$ws.Range("K132").Value = 13275
$ws.Range("M132").Value = -10745

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 215.6
$ws.Range("I2").Value = 99.666664
$ws.Range("J2").Value = 389.5
$ws.Range("K2").Value = 99.666664
$ws.Range("L2").Value = 389.5
$ws.Range("M2").Value = 13.333336
$ws.Range("N2").Value = -615.5
$ws.Range("H107").Value = 1436.6666
$ws.Range("I107").Value = 573.8333
$ws.Range("K107").Value = 573.8333
$ws.Range("M107").Value = 1346.1667
$ws.Range("H113").Value = 1830.25
$ws.Range("I113").Value = 1724.8182
$ws.Range("K113").Value = 1724.8182
$ws.Range("M113").Value = 445.1818000000001
$ws.Range("H132").Value = 3748.7317
$ws.Range("J132").Value = 5065.467
$ws.Range("L132").Value = 15196.401
$ws.Range("N132").Value = -20256.401

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2405.5454
$ws.Range("I46").Value = 989.3333
$ws.Range("J46").Value = 2936.625
$ws.Range("K46").Value = 989.3333
$ws.Range("L46").Value = 2936.625
$ws.Range("M46").Value = -801.3333
$ws.Range("N46").Value = -3312.625
$ws.Range("H55").Value = 356.55
$ws.Range("I55").Value = 349.35294
$ws.Range("J55").Value = 397.33334
$ws.Range("K55").Value = 349.35294
$ws.Range("L55").Value = 397.33334
$ws.Range("M55").Value = -176.35294
$ws.Range("N55").Value = -743.33334
$ws.Range("H132").Value = 4758.05
$ws.Range("I132").Value = 2859.1365
$ws.Range("K132").Value = 8577.4095
$ws.Range("M132").Value = -6047.4095
$ws.Range("H136").Value = 2491.25
$ws.Range("I136").Value = 2246.8823
$ws.Range("K136").Value = 6740.646900000001
$ws.Range("M136").Value = -4190.646900000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H122").Value = 1833.3226
$ws.Range("I122").Value = 1830.8518
$ws.Range("K122").Value = 5492.555399999999
$ws.Range("M122").Value = -3042.555399999999
$ws.Range("H132").Value = 4100.606
$ws.Range("I132").Value = 3604.4814
$ws.Range("K132").Value = 10813.4442
$ws.Range("M132").Value = -8283.4442
$ws.Range("H136").Value = 3849.7856
$ws.Range("I136").Value = 1650.8
$ws.Range("K136").Value = 3849.7856
$ws.Range("M136").Value = -2402.4
